# Automatische test-sync: 2025-08-03 18:55:50
# Appends a new incoming-mail log row (#20) to the "Logs" sheet and
# refreshes the "Dashboard" category-count summary to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Logs" sheet: append row 48 with the new test mail.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A48").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("B48").Value = "mailmind.test@zohomail.eu"
$logs.Range("C48").Value = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D48").Value = "Klacht / Probleem"
$logs.Range("E48").Value = "Bedankt, we hebben dit doorgestuurd naar klachten@bedrijf.nl."
$logs.Range("F48").Value = "2025-08-03 18:55:27"
$logs.Range("G48").Value = "Ja"
$logs.Range("H48").Value = "Ja"
$logs.Range("I48").Value = "Nee"
$logs.Range("J48").Value = "Nee"

# Extend the conditional-formatting ranges (D/G/H/I/J) from row 47 down
# to the newly added row 48, matching the grown data range.
$logs.Range("D2:D47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D48"))
$logs.Range("G2:G47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G48"))
$logs.Range("H2:H47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H48"))
$logs.Range("I2:I47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I48"))
$logs.Range("J2:J47").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J48"))

# ---------------------------------------------------------------------
# 2) "Dashboard" sheet: "Klacht / Probleem" now has 2 occurrences, so it
#    moves up in the (count desc) ranking, shuffling rows 7-9.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Klacht / Probleem"
$dash.Range("A8").Value = "Klantenservice / Contact"
$dash.Range("A9").Value = "Retour / Terugbetaling"
$dash.Range("B9").Value = 2
